$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-05-26"

# Update the label in column A row 6 (May through date)
$ws.Range("A6").Value = "May (through 05-26)"

# Update the May row (row 6) values
$ws.Range("B6").Value = 16
$ws.Range("C6").Value = 39
$ws.Range("D6").Value = 51
$ws.Range("E6").Value = 42
$ws.Range("F6").Value = 37
$ws.Range("G6").Value = 51
$ws.Range("H6").Value = 98
$ws.Range("I6").Value = 95

# Update the Total row (row 7) values
$ws.Range("B7").Value = 105
$ws.Range("C7").Value = 201
$ws.Range("D7").Value = 304
$ws.Range("E7").Value = 288
$ws.Range("F7").Value = 192
$ws.Range("G7").Value = 313
$ws.Range("H7").Value = 621
$ws.Range("I7").Value = 646
